# Sprint 3 rubric update — User Story #23 completed (30 minutes), and
# "the owner is now tied to all rental items" items get their points.
#
# Rubric rows affected (Points column, C):
#   Row 8  - "Add rentals to catalog, shopping cart, and checkout"  -> 10
#   Row 16 - "Include products that artisans offer for sale"        -> 10
#   Row 20 - "Handle credit cards with REST calls (requests library)" -> 15
#
# All of the SUM()/totals/percentage formulas recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 10
$ws.Range("C16").Value = 10
$ws.Range("C20").Value = 15

# Reflect the author's final cursor position / on-screen selection.
[void]$ws.Range("C9").Select()
